$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "SVM"

$ws.Range("A2").Value = "Name"
$ws.Range("B2").Value = "linear"
$ws.Range("C2").Value = "poly"
$ws.Range("D2").Value = "rbf"
$ws.Range("E2").Value = "sigmoid"
$ws.Range("F2").Value = "precomputed"

$ws.Range("A3").Value = "train score"
$ws.Range("B3").Value = 0.9776
$ws.Range("C3").Value = 0.9776
$ws.Range("D3").Value = 0.9805
$ws.Range("E3").Value = 0.8263

$ws.Range("A4").Value = "test score"
$ws.Range("B4").Value = 0.9779
$ws.Range("C4").Value = 0.9771
$ws.Range("D4").Value = 0.9803
$ws.Range("E4").Value = 0.8279

$ws.Range("F3").Value = "N/A"
$ws.Range("F4").Value = "N/A"

$ws.Columns.Item(1).ColumnWidth = 10.66
$ws.Columns.Item(6).ColumnWidth = 13.66
$ws.Columns.Item(7).ColumnWidth = 15.0

$ws.Range("A8").Select()
